# ecs事件 增加beforeupdate 和 lateupdate
#
# The "描述" (description) column — column B, containing the repeated "xx"
# placeholder value — is removed from every sheet (Sheet1/Sheet2/Sheet3),
# shifting the 类型/中文/英语 (type/Chinese/English) columns one slot to
# the left. Outline level metadata is trimmed to match (one fewer grouped
# column), and the selections left behind mirror the edit session: Sheet2
# had column B selected (about to be deleted), Sheet1 ends up the active
# tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Remove the description column (B) from each sheet; cells to the right
# shift left to fill the gap.
$ws1.Columns("B").Delete()
$ws1.Outline.ShowLevels(5, 3)

$ws2.Columns("B").Delete()
$ws2.Outline.ShowLevels(5, 3)

$ws3.Columns("B").Delete()
$ws3.Outline.ShowLevels(5, 3)

# Leave the selection/active-sheet state matching the end of the edit
# session: Sheet3 and Sheet2 were visited (Sheet2 with the whole of the
# just-deleted column B still marked), finishing back on Sheet1.
$ws3.Select()
$ws3.Range("G21").Select()

$ws2.Select()
$ws2.Range("B1:B1048576").Select()

$ws1.Select()
$ws1.Range("G13").Select()
